# Scheduled-runner update: refresh computed profit/price columns (H:N) for a
# set of leve rows across the ALC/ARM/CRP/CUL/GSM/LTW/WVR sheets.
# Columns: H=currentAveragePrice, I=currentAveragePriceNQ, J=currentAveragePriceHQ,
#          K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ, N=LeveProfitHQ
$wb = $excel.ActiveWorkbook

# ALC row 17
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2631.08
$ws.Range("I17").Value = 997
$ws.Range("J17").Value = 2773.1738
$ws.Range("K17").Value = 2991
$ws.Range("L17").Value = 8319.5214
$ws.Range("M17").Value = -2823
$ws.Range("N17").Value = -8655.5214

# ALC row 28
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 7207
$ws.Range("I28").Value = 1590
$ws.Range("K28").Value = 1590
$ws.Range("M28").Value = -1105

# ALC row 33
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 133.5
$ws.Range("I33").Value = 94
$ws.Range("K33").Value = 94
$ws.Range("M33").Value = 135

# ALC row 94
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H94").Value = 3523.0908
$ws.Range("I94").Value = 3523.0908
$ws.Range("K94").Value = 3523.0908
$ws.Range("M94").Value = -3072.0908

# ALC row 97
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H97").Value = 1000
$ws.Range("J97").Value = 1000
$ws.Range("L97").Value = 3000
$ws.Range("N97").Value = -3992

# ALC row 100
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 2327.8696
$ws.Range("J100").Value = 2537.182
$ws.Range("L100").Value = 2537.182
$ws.Range("N100").Value = -3619.182

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1929.7916
$ws.Range("I137").Value = 1033.5385
$ws.Range("K137").Value = 3100.6155
$ws.Range("M137").Value = -550.6155000000003

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 6252.8
$ws.Range("I74").Value = 6141.4165
$ws.Range("K74").Value = 6141.4165
$ws.Range("M74").Value = -5267.4165

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 6252.8
$ws.Range("I77").Value = 6141.4165
$ws.Range("K77").Value = 30707.0825
$ws.Range("M77").Value = -26339.0825

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6047.697
$ws.Range("I31").Value = 3724.75
$ws.Range("K31").Value = 3724.75
$ws.Range("M31").Value = -3429.75

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 6047.697
$ws.Range("I34").Value = 3724.75
$ws.Range("K34").Value = 3724.75
$ws.Range("M34").Value = -3522.75

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3394.4
$ws.Range("I58").Value = 1661.6666
$ws.Range("J58").Value = 5993.5
$ws.Range("K58").Value = 1661.6666
$ws.Range("L58").Value = 5993.5
$ws.Range("M58").Value = -1458.6666
$ws.Range("N58").Value = -6399.5

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 3394.4
$ws.Range("I136").Value = 1661.6666
$ws.Range("J136").Value = 5993.5
$ws.Range("K136").Value = 4984.9998
$ws.Range("L136").Value = 17980.5
$ws.Range("M136").Value = -2434.9998
$ws.Range("N136").Value = -23080.5

# CRP row 141
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 216099.5
$ws.Range("J141").Value = 216099.5
$ws.Range("L141").Value = 216099.5
$ws.Range("N141").Value = -226459.5

# CUL row 44
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 1166.5
$ws.Range("I44").Value = 500
$ws.Range("K44").Value = 1500
$ws.Range("M44").Value = -1102

# CUL row 55
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 2856.1428
$ws.Range("J55").Value = 3599.7778
$ws.Range("L55").Value = 10799.3334
$ws.Range("N55").Value = -11153.3334

# CUL row 57
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H57").Value = 2250
$ws.Range("J57").Value = 2562.5
$ws.Range("L57").Value = 7687.5
$ws.Range("N57").Value = -8805.5

# CUL row 60
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 1027.8182
$ws.Range("J60").Value = 1763.4546
$ws.Range("L60").Value = 5290.3638
$ws.Range("N60").Value = -5792.3638

# CUL row 109
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 1239.3077
$ws.Range("I109").Value = 940.1111
$ws.Range("K109").Value = 2820.3333
$ws.Range("M109").Value = -1780.3333

# CUL row 118
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H118").Value = 3682.3333
$ws.Range("I118").Value = 3682.3333
$ws.Range("K118").Value = 11046.9999
$ws.Range("M118").Value = -9803.999899999999

# CUL row 121
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 34000
$ws.Range("J121").Value = 1000
$ws.Range("L121").Value = 3000
$ws.Range("N121").Value = -5620

# CUL row 122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 826.9
$ws.Range("I122").Value = 879.5
$ws.Range("J122").Value = 748
$ws.Range("K122").Value = 7915.5
$ws.Range("L122").Value = 6732
$ws.Range("M122").Value = -5465.5
$ws.Range("N122").Value = -11632

# CUL row 139
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 1665.8667
$ws.Range("I139").Value = 1077.9231
$ws.Range("J139").Value = 5487.5
$ws.Range("K139").Value = 3233.7693
$ws.Range("L139").Value = 16462.5
$ws.Range("M139").Value = 1906.2307
$ws.Range("N139").Value = -26742.5

# GSM row 92
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 4048.6
$ws.Range("J92").Value = 4048.6
$ws.Range("L92").Value = 4048.6
$ws.Range("N92").Value = -7792.6

# GSM row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2816.3333
$ws.Range("I126").Value = 1724.5
$ws.Range("K126").Value = 5173.5
$ws.Range("M126").Value = -2703.5

# LTW row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3805.5
$ws.Range("I7").Value = 3805.5
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 3805.5
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -3693.5
$ws.Range("N7").ClearContents()

# LTW row 32
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 8519.714
$ws.Range("I32").Value = 6273
$ws.Range("K32").Value = 6273
$ws.Range("M32").Value = -5956

# LTW row 68
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 4560.4
$ws.Range("I68").Value = 4560.4
$ws.Range("K68").Value = 4560.4
$ws.Range("M68").Value = -3811.4

# LTW row 71
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 4560.4
$ws.Range("I71").Value = 4560.4
$ws.Range("K71").Value = 22802
$ws.Range("M71").Value = -19058

# LTW row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 3805.5
$ws.Range("I126").Value = 3805.5
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 11416.5
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -8946.5
$ws.Range("N126").ClearContents()

# WVR row 31
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 1500
$ws.Range("I31").Value = 1500
$ws.Range("K31").Value = 1500
$ws.Range("M31").Value = -1152

# WVR row 43
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 56686.668
$ws.Range("I43").Value = 40000
$ws.Range("J43").Value = 65030
$ws.Range("K43").Value = 40000
$ws.Range("L43").Value = 65030
$ws.Range("M43").Value = -39851
$ws.Range("N43").Value = -65328

# WVR row 100
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1874.7142
$ws.Range("I100").Value = 1874.7142
$ws.Range("K100").Value = 3749.4284
$ws.Range("M100").Value = -3208.4284

# WVR row 135
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H135").Value = 80000
$ws.Range("J135").Value = 80000
$ws.Range("L135").Value = 80000
$ws.Range("N135").Value = -90140
